# WIP update:
#  - rename "Sheet3" -> "Parser failures" (defined name "failed_1" auto-updates
#    its formula reference since it points at this sheet)
#  - scroll the "latest" sheet view down so row 21 is the first visible row
#  - make "Parser failures" the active/selected tab (it was "unmapped-latest")

$wb = $excel.ActiveWorkbook

# Scroll "latest" so its viewport starts at A21 (no freeze/split involved).
$wsLatest = $wb.Worksheets.Item("latest")
$wsLatest.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1

# Rename Sheet3 -> "Parser failures" (updates the failed_1 defined name too).
$wsParserFailures = $wb.Worksheets.Item("Sheet3")
$wsParserFailures.Name = "Parser failures"

# Make it the active tab (moves tabSelected off "unmapped-latest").
$wsParserFailures.Activate()
